# framework.xlsx - "Further changes to framework, js-files"
#
# Adds a new "MIF" / "MIF_VISIT" pair of forms:
#  - choices sheet: two new choice rows for the `forms` list
#  - survey sheet: two new branch blocks (label / external_link / exit section)
#  - settings sheet: bump form_version
#  - leaves the final selection on the survey sheet, matching the source file.

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("settings")
$wsChoices  = $wb.Worksheets.Item("choices")
$wsSurvey   = $wb.Worksheets.Item("survey")
$wsFrameworkTranslations = $wb.Worksheets.Item("framework_translations")

# ---------------------------------------------------------------------------
# settings: bump form_version
# ---------------------------------------------------------------------------
$wsSettings.Range("B4").Value = 20230729

# ---------------------------------------------------------------------------
# choices: register the two new "forms" choice-list entries (MIF / MIF_VISIT)
# Order of assignment matters: it controls the order new strings land in the
# shared-string table (MIF, MIF_VISIT, Core MIF, MIF visit, MIF visita, ...).
# ---------------------------------------------------------------------------
$wsChoices.Range("B6").Value = "MIF"
$wsChoices.Range("B7").Value = "MIF_VISIT"
$wsChoices.Range("C6").Value = "Core MIF"
$wsChoices.Range("D6").Value = "Core MIF"
$wsChoices.Range("C7").Value = "MIF visit"
$wsChoices.Range("D7").Value = "MIF visita"
$wsChoices.Range("A6").Value = "forms"
$wsChoices.Range("A7").Value = "forms"

# ---------------------------------------------------------------------------
# survey: add the two new branch blocks, mirroring the existing
# CHILDREN/CHILDFU/PREGNANCIES/PREGNANCYFU pattern (rows 8-19).
# ---------------------------------------------------------------------------

# Clear out the two stray formatted-but-empty placeholder cells first.
$wsSurvey.Range("B20").Clear()
$wsSurvey.Range("B23").Clear()

# -- MIF block --
$wsSurvey.Range("A20").Value = "MIF"
$wsSurvey.Range("B21").Value = "''?' + odkSurvey.getHashString('MIF')"
$wsSurvey.Range("E21").Value = "external_link"
$wsSurvey.Range("G21").Value = "Open form"
$wsSurvey.Range("C22").Value = "exit section"

# -- MIF_VISIT block --
$wsSurvey.Range("A23").Value = "MIF_VISIT"
$wsSurvey.Range("B24").Value = "''?' + odkSurvey.getHashString('MIF_VISIT')"
$wsSurvey.Range("E24").Value = "external_link"
$wsSurvey.Range("G24").Value = "Open form"
$wsSurvey.Range("C25").Value = "exit section"

# Re-apply the quote-prefix style used by the sibling formula cells
# (B9/B12/B15/B18) to the two new formula cells.
$wsSurvey.Range("B21").Style = $wsSurvey.Range("B18").Style
$wsSurvey.Range("B24").Style = $wsSurvey.Range("B18").Style

# ---------------------------------------------------------------------------
# Final view state: survey sheet ends up active, with B25 selected.
# ---------------------------------------------------------------------------
$wsFrameworkTranslations.Activate()
$wsFrameworkTranslations.Range("B36").Select()

$wsChoices.Activate()
$wsChoices.Range("C10").Select()

$wsSurvey.Activate()
$wsSurvey.Range("B25").Select()
